# Generate Report for Handback
# Refreshes the handback-status workbook for a newly generated run:
#   - the two source markdown files got new GUID-based names
#   - the zh-cn / de-de xliff correspondence for the second file now
#     reuses the (re-generated) artifacts of the first file
#   - handoff / handback / xliff-generate timestamps move forward

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldGuid1 = "562117bb-5af2-40a6-8713-b398be52c1ed"
$newGuid1 = "462b48b3-567d-4863-b338-db92d04a245e"

$oldGuid2 = "5eb3d5b0-c04d-4b6b-83f1-cd1ac62f5fc6"
$newGuid2 = "ffffb957a9e2-d226-4961-a859-d90d0a82d103"

$newXlfZhCn = "$newGuid1.afeafe056b3d5c2feb2964de8ab7465eacddbd50.zh-cn.xlf"
$newXlfDeDe = "$newGuid1.afeafe056b3d5c2feb2964de8ab7465eacddbd50.de-de.xlf"

# --- Overview sheet -------------------------------------------------
$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-09-06 05:19:30"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-09-06 05:19:30"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = $newXlfZhCn
$wsZhCn.Range("H2").Value = "2016-09-06 05:19:25"
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = $newXlfZhCn
$wsZhCn.Range("K2").Value = "2016-09-06 05:19:43"

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = $newXlfZhCn
$wsZhCn.Range("H3").Value = "2016-09-06 05:19:25"
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = $newXlfZhCn
$wsZhCn.Range("K3").Value = "2016-09-06 05:19:43"

# --- de-de sheet ------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = $newXlfDeDe
$wsDeDe.Range("H2").Value = "2016-09-06 05:19:30"
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = $newXlfDeDe
$wsDeDe.Range("K2").Value = "2016-09-06 05:19:51"

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = $newXlfDeDe
$wsDeDe.Range("H3").Value = "2016-09-06 05:19:30"
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = $newXlfDeDe
$wsDeDe.Range("K3").Value = "2016-09-06 05:19:51"

# --- Hyperlink display text --------------------------------------------
# The runtime's Hyperlinks.Delete() clears every hyperlink on the sheet
# (not just the targeted range), and reading back .Address afterwards
# isn't reliable, so rebuild each sheet's hyperlinks in one pass: clear
# them all, then re-add them in the same order using their original
# target URLs (unchanged by this edit) with the refreshed display text.
# Re-adding with identical target URLs in the same order reproduces the
# original rId assignments in the worksheet's relationship part.

$urlGuid1Base = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e33cb351585ead2ebb1cd6a3396bac73f18e5e8/e2e/$oldGuid1.md"
$urlGuid2Base = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e33cb351585ead2ebb1cd6a3396bac73f18e5e8/e2e/$oldGuid2.md"
$urlGuid1ZhCn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e5b9278ed750ddea3dbbd35883cc974803bff2fc/e2e/$oldGuid1.md"
$urlGuid2ZhCn = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e5b9278ed750ddea3dbbd35883cc974803bff2fc/e2e/$oldGuid2.md"
$urlGuid1DeDe = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/356f5adc18aeacec45ef4692f2a564520ba1726c/e2e/$oldGuid1.md"
$urlGuid2DeDe = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/356f5adc18aeacec45ef4692f2a564520ba1726c/e2e/$oldGuid2.md"

# --- Overview sheet hyperlinks (B2, B3) ---
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $urlGuid1Base, "", "", "e2e\$newGuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $urlGuid2Base, "", "", "e2e\$newGuid2.md")

# --- zh-cn sheet hyperlinks (A2, I2, A3, I3) ---
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlGuid1Base, "", "", "$newGuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlGuid1ZhCn, "", "", "$newGuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlGuid2Base, "", "", "$newGuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlGuid2ZhCn, "", "", "$newGuid2.md")

# --- de-de sheet hyperlinks (A2, I2, A3, I3) ---
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlGuid1Base, "", "", "$newGuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlGuid1DeDe, "", "", "$newGuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlGuid2Base, "", "", "$newGuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlGuid2DeDe, "", "", "$newGuid2.md")
